$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-04-22 Tuesday" "2025-04-23 Wednesday"

Replace-Text "45÷5=" "84÷3="
Replace-Text "48÷7=" "23÷4="
Replace-Text "86÷8=" "11÷4="
Replace-Text "38÷9=" "52÷2="
Replace-Text "62÷9=" "43÷6="

Replace-Text "71÷2=" "22÷6="
Replace-Text "83÷7=" "19÷2="
Replace-Text "11÷2=" "78÷2="
Replace-Text "35÷4=" "47÷3="
Replace-Text "92÷9=" "63÷9="

Replace-Text "10÷8=" "61÷9="
Replace-Text "72÷2=" "90÷4="
Replace-Text "32÷3=" "41÷3="
Replace-Text "99÷7=" "23÷2="
Replace-Text "67÷5=" "54÷6="

Replace-Text "73÷3=" "95÷9="
Replace-Text "68÷6=" "97÷9="
Replace-Text "57÷2=" "76÷7="
Replace-Text "38÷5=" "34÷7="
Replace-Text "42÷4=" "74÷6="

Replace-Text "19÷9=" "47÷8="
Replace-Text "37÷5=" "14÷3="
Replace-Text "59÷7=" "75÷8="
Replace-Text "15÷8=" "99÷8="
Replace-Text "28÷8=" "74÷4="
